$wb = $excel.ActiveWorkbook

# --- 1. Split the existing "2022-Q3" sheet into "2022-Q4" (new data) + a
#        fresh "2022-Q3" sheet that keeps the original Q3 data. ---

$q3Sheet = $wb.Worksheets.Item("2022-Q3")

# Duplicate it first so the copy retains all the existing Q3 values/styles;
# place the copy right after the original.
$q3Sheet.Copy($null, $q3Sheet)
$newQ3Sheet = $wb.Worksheets.Item("2022-Q3 (2)")

# The original sheet becomes "2022-Q4" and gets overwritten with the new
# quarter's figures. Rename it first so the freshly-copied sheet can then
# take over the "2022-Q3" name.
$q4Sheet = $q3Sheet
$q4Sheet.Name = "2022-Q4"
$newQ3Sheet.Name = "2022-Q3"

# --- 2. Overwrite the (now) "2022-Q4" sheet with the Q4 fund data. ---

$q4Sheet.Range("A2").Value = 0
$q4Sheet.Range("B2").Value = "'007811"
$q4Sheet.Range("C2").Value = "淳厚信泽灵活配置混合A"
$q4Sheet.Range("D2").Value = "'4.32"
$q4Sheet.Range("E2").Value = "'79.09"
$q4Sheet.Range("F2").Value = "'2.78"
$q4Sheet.Range("G2").Value = "'0.1201"
$q4Sheet.Range("H2").Value = 10

$q4Sheet.Range("A3").Value = 1
$q4Sheet.Range("B3").Value = "'010551"
$q4Sheet.Range("C3").Value = "淳厚欣颐一年持有期混合"
$q4Sheet.Range("D3").Value = "'2.49"
$q4Sheet.Range("E3").Value = "'88.75"
$q4Sheet.Range("F3").Value = "'2.77"
$q4Sheet.Range("G3").Value = "'0.0690"
$q4Sheet.Range("H3").Value = 10

$q4Sheet.Range("A4").Value = 2
$q4Sheet.Range("B4").Value = "'011349"
$q4Sheet.Range("C4").Value = "淳厚现代服务业股票A"
$q4Sheet.Range("D4").Value = "'2.41"
$q4Sheet.Range("E4").Value = "'88.44"
$q4Sheet.Range("F4").Value = "'2.76"
$q4Sheet.Range("G4").Value = "'0.0665"
$q4Sheet.Range("H4").Value = 10

$q4Sheet.Range("A5").Value = 3
$q4Sheet.Range("B5").Value = "'007812"
$q4Sheet.Range("C5").Value = "淳厚信泽灵活配置混合C"
$q4Sheet.Range("D5").Value = "'0.81"
$q4Sheet.Range("E5").Value = "'79.09"
$q4Sheet.Range("F5").Value = "'2.78"
$q4Sheet.Range("G5").Value = "'0.0225"
$q4Sheet.Range("H5").Value = 10

$q4Sheet.Range("A6").Value = 4
$q4Sheet.Range("B6").Value = "'011350"
$q4Sheet.Range("C6").Value = "淳厚现代服务业股票C"
$q4Sheet.Range("D6").Value = "'0.57"
$q4Sheet.Range("E6").Value = "'88.44"
$q4Sheet.Range("F6").Value = "'2.76"
$q4Sheet.Range("G6").Value = "'0.0157"
$q4Sheet.Range("H6").Value = 10

# --- 3. Update the "总计" (totals) summary sheet. ---

$totalSheet = $wb.Worksheets.Item("总计")

# Existing row 2 now describes 2022-Q4.
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("D2").Value = "'0.29"

# New row 3 carries the original 2022-Q3 totals.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q3"
$totalSheet.Range("C3").Value = 5
$totalSheet.Range("D3").Value = "'0.25"
